$d = $word.ActiveDocument

# Move to the very end of the document body (after the last paragraph mark).
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd = 0

# Paragraph 1: empty
$end.InsertParagraphAfter()
$end.Collapse(0)

# Paragraph 2: empty
$end.InsertParagraphAfter()
$end.Collapse(0)

# Paragraph 3: empty
$end.InsertParagraphAfter()
$end.Collapse(0)

# Paragraph 4: "Cinemática"
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "Cinemática"

$end = $d.Content
$end.Collapse(0)

# Paragraph 5: "Cambiar textos"
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "Cambiar textos"

$end = $d.Content
$end.Collapse(0)

# Paragraph 6: "Cambiar musica"
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "Cambiar musica"
